# issue #5: stock data from json to db
# Adds "category", "source_file" and "index" columns to the 股票 (stock)
# sheet, matching the new normalized-export schema.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at I (pushes the old date/legislator_name/legislator_id
# columns from I/J/K to J/K/L) and two new columns at the end (M, N).
$ws.Columns("I:I").Insert()
$ws.Columns("M:N").Insert()

# Header row
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data rows: category = "normal", source_file = "tmp2bc41",
# index = same value as column A (the row's record index).
for ($r = 2; $r -le 13; $r++) {
    $recordIndex = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp2bc41"
    $ws.Cells.Item($r, 14).Value = $recordIndex
}
